$wb = $excel.ActiveWorkbook

# --- Sheet "ROW50-FE-LIFTER": append new row 55 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 55
$ws1.Cells.Item($r, 1).Value = 45751.70443925926
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item($r - 1, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x5e"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = 568631262647113970876416.0
$ws1.Cells.Item($r, 8).Value = 350
$ws1.Cells.Item($r, 9).Value = 14

# --- Sheet "ROW50-MID-LIFTER": append new row 57 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 57
$ws2.Cells.Item($r, 1).Value = 45751.67511574074
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item($r - 1, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x66"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400

# G57 must be stored as TEXT (matches the rest of this sheet's "ID_DEC"
# column, which is all inlineStr) even though it looks like a huge
# integer. Plain `.Value = "<digits>"` gets auto-coerced to a Number by
# Excel's smart typing. Instead, stage a formula in a scratch cell whose
# *result* is that text string, copy it, and paste-special *values only*
# into G57 - only the literal text lands in the cell, with no number
# format / quote-prefix style residue; then wipe the scratch cell.
$scratch = $ws2.Cells.Item(200, 1)
$scratch.Formula = "=""568631262647113771663628"""
$scratch.Copy()
$ws2.Cells.Item($r, 7).PasteSpecial(-4163)
$scratch.Clear()

$ws2.Cells.Item($r, 8).Value = 358
$ws2.Cells.Item($r, 9).Value = 25

# --- Sheet "ROW11-FE-LIFTER": append new row 55 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 55
$ws3.Cells.Item($r, 1).Value = 45751.73562108797
$ws3.Cells.Item($r, 1).NumberFormat = $ws3.Cells.Item($r - 1, 1).NumberFormat
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x5e"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = 568631262647113970876416.0
$ws3.Cells.Item($r, 8).Value = 350
$ws3.Cells.Item($r, 9).Value = 20

# --- Sheet "ROW11-MID-LIFTER": append new row 55 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 55
$ws4.Cells.Item($r, 1).Value = 45751.86985818287
$ws4.Cells.Item($r, 1).NumberFormat = $ws4.Cells.Item($r - 1, 1).NumberFormat
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x66"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = 568631262647113970876416.0
$ws4.Cells.Item($r, 8).Value = 358
$ws4.Cells.Item($r, 9).Value = 25
